# Update the "Agenda" sheet: move the fixed "Hours Required" figures for
# AI Concepts/DevOps/Python/Quantum Physics/Block Chain/Pre-Calculus out of
# row 10 into a brand-new row 19 (so row 10 reflects only the newly-added
# "Fast Reading" hours), and bump the Fast Reading total (column P) up to 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")

# Preserve the existing "Hours Required" row formatting (s="61") on the new
# row by copying it down before overwriting row 10's values.
$ws.Range("C10:H10").Copy()
$ws.Range("C19:H19").PasteSpecial(-4122)

# Re-home the old hours-required numbers on row 19.
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 7

# Row 10 (Hours Required) now has no requirement for those subjects, and the
# Fast Reading (column P) requirement grows from 3.5 to 16 hours.
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("P10").Value = 16

# Match the author's last on-sheet selection.
$ws.Activate()
$ws.Range("I14").Select()
